# Update cryptocurrency price/volume data (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.520.24"
$ws.Range("E2").Value = "  -4.45%  "
$ws.Range("D3").Value = "3.008.75"
$ws.Range("E3").Value = "  -5.70%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.96"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.34"
$ws.Range("E6").Value = "  -6.64%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "3.005.24"
$ws.Range("E8").Value = "  -5.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("E10").Value = "  -6.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.16"
$ws.Range("E11").Value = "  -2.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -5.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.69"
$ws.Range("E14").Value = "  -5.79%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "3.498.09"
$ws.Range("E16").Value = "  -5.84%  "
$ws.Range("D17").Value = "3.007.23"
$ws.Range("E17").Value = "  -5.57%  "
$ws.Range("D18").Value = "60.359.69"
$ws.Range("E18").Value = "  -4.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.56"
$ws.Range("E20").Value = "  -6.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.18"
$ws.Range("E21").Value = "  -5.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.665"
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.06"
$ws.Range("E23").Value = "  -7.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.94"
$ws.Range("E24").Value = "  -3.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.50"
$ws.Range("E25").Value = "  -3.86%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.57"
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  -5.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.91"
$ws.Range("E30").Value = "  -6.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.20"
$ws.Range("E31").Value = "  -8.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.36"
$ws.Range("E32").Value = "  -7.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0943"
$ws.Range("E33").Value = "  -8.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("E34").Value = "  -9.03%  "
$ws.Range("E35").Value = "  -7.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.64"
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.11"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("D38").Value = "0.0₃0668"
$ws.Range("E38").Value = "  -7.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.43"
$ws.Range("E39").Value = "  +3.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0361"
$ws.Range("E40").Value = "  -7.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "384.27"
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.50"
$ws.Range("E43").Value = "  -7.74%  "
$ws.Range("D44").Value = "2.658.94"
$ws.Range("E44").Value = "  -5.57%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.237"
$ws.Range("E46").Value = "  -6.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.03"
$ws.Range("E47").Value = "  -5.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.16"
$ws.Range("E48").Value = "  -7.13%  "
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.88"
$ws.Range("E50").Value = "  -6.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.136"
$ws.Range("E51").Value = "  +4.89%  "
